# correction to Asn flux
# The "Asparagine pos" column (column F) is removed entirely (header +
# data), which shifts every later column one position to the left
# (Asparagine-13C4 pos lands in the old G -> new F slot, etc.).
# The values that land in the corrected "Asparagine-13C4 pos" column
# (new column F) are also replaced with the corrected flux numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "Asparagine pos" column (old column F); this shifts
# columns G:U left to F:T and renumbers the shared strings automatically.
$ws.Columns("F").Delete()

# Write the corrected "Asparagine-13C4 pos" flux values (now in column F)
# for the data rows that previously held numeric flux results.
$ws.Range("F14").Value = 3.038633187518536
$ws.Range("F15").Value = 2.243094108572594
$ws.Range("F16").Value = 2.068234961241354
$ws.Range("F17").Value = 2.976720607509004
$ws.Range("F18").Value = 2.225486338339457
$ws.Range("F19").Value = 2.588102911237619
$ws.Range("F20").Value = 2.268791010543967
$ws.Range("F21").Value = 2.34093061538717
$ws.Range("F22").Value = 2.137606490472667
$ws.Range("F23").Value = 1.955318440522744
$ws.Range("F24").Value = 2.057154078031326
$ws.Range("F25").Value = 2.142724429855566
